$d = $word.ActiveDocument

# Returns the 1-based index of the first paragraph (at/after $startIndex)
# whose text equals $text (ignoring the trailing paragraph mark).
function Find-ParagraphIndex($doc, $text, $startIndex) {
    for ($i = $startIndex; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.TrimEnd("`r") -eq $text) {
            return $i
        }
    }
    return -1
}

# Anchor on the "Commande : PLAY" section so every subsequent lookup is
# scoped to this command's block (several other commands share line text
# such as "Réponses :").
$playCmdIdx = Find-ParagraphIndex $d "Commande : PLAY" 1

# --- 1) "Paramètres : <index>" -> "Paramètres : <index> <color>" -----------
$paramIdx = Find-ParagraphIndex $d "Paramètres : <index>" $playCmdIdx
$paramRange = $d.Paragraphs($paramIdx).Range
$paramRange.Find.Execute("Paramètres : <index>", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Paramètres : <index> <color>", 2)

# --- 2) New paragraph right after the PLAY description, right before
#        "Réponses :", explaining the new <color> parameter ----------------
$descIdx = Find-ParagraphIndex $d "Permets de jouer une carte si c’est votre tour." $playCmdIdx
$d.Paragraphs($descIdx).Range.InsertParagraphAfter()
$d.Paragraphs($descIdx + 1).Range.Text = `
    "Le paramètre <color> est indispensable si vous jouez un ChangeColor ou un Plus4"

# --- 3) Fix the BADINDEX error message wording (inexistant -> inexistante) -
$badIdx = Find-ParagraphIndex $d "424 ERR_BADINDEX : type ou couleur de carte inexistant" $playCmdIdx
$badRange = $d.Paragraphs($badIdx).Range
$badRange.Find.Execute("424 ERR_BADINDEX : type ou couleur de carte inexistant", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "424 ERR_BADINDEX : type ou couleur de carte inexistante", 2)

# --- 4) New paragraph right after it introducing the 425 UNDEFINEDCOLOR
#        error response, right before "200 PLAY_OK : jeu confirmé" ---------
$d.Paragraphs($badIdx).Range.InsertParagraphAfter()
$d.Paragraphs($badIdx + 1).Range.Text = `
    "425 ERR_UNDEFINEDCOLOR : couleur de carte inexistante"
